$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the last used row in column C (data starts at row 2)
$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row

# Column C ("Förändrad") holds a date serial that was bumped by 2 days
# (45175 -> 45177) for every data row.
$rng = $ws.Range("C2:C$lastRow")
$rng.Value = 45177
